$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1731.3158
$ws.Range("I107").Value = 1506.8462
$ws.Range("K107").Value = 1506.8462
$ws.Range("M107").Value = 413.1538
$ws.Range("H137").Value = 1108.4814
$ws.Range("I137").Value = 1098.909
$ws.Range("J137").Value = 1150.6
$ws.Range("K137").Value = 3296.727
$ws.Range("L137").Value = 3451.8
$ws.Range("M137").Value = -746.7270000000003
$ws.Range("N137").Value = -8551.799999999999
$ws.Range("H138").Value = 1364.4366
$ws.Range("I138").Value = 738.0357
$ws.Range("J138").Value = 1772.3256
$ws.Range("K138").Value = 2214.1071
$ws.Range("L138").Value = 5316.976799999999
$ws.Range("M138").Value = 2925.8929
$ws.Range("N138").Value = -15596.9768

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7598.7334
$ws.Range("I2").Value = 912.9167
$ws.Range("J2").Value = 34342
$ws.Range("K2").Value = 912.9167
$ws.Range("L2").Value = 34342
$ws.Range("M2").Value = -799.9167
$ws.Range("N2").Value = -34568
$ws.Range("H32").Value = 4175.5884
$ws.Range("I32").Value = 4649.6553
$ws.Range("K32").Value = 4649.6553
$ws.Range("M32").Value = -4362.6553
$ws.Range("H97").Value = 511.125
$ws.Range("I97").Value = 369.85715
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 369.85715
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 126.14285
$ws.Range("N97").Value = -2492
$ws.Range("H106").Value = 22370
$ws.Range("J106").Value = 22370
$ws.Range("L106").Value = 22370
$ws.Range("N106").Value = -24894
$ws.Range("H116").Value = 7598.7334
$ws.Range("I116").Value = 912.9167
$ws.Range("J116").Value = 34342
$ws.Range("K116").Value = 912.9167
$ws.Range("L116").Value = 34342
$ws.Range("M116").Value = 1381.0833
$ws.Range("N116").Value = -38930
$ws.Range("H122").Value = 1540.5555
$ws.Range("I122").Value = 1414.4667
$ws.Range("J122").Value = 2171
$ws.Range("K122").Value = 4243.4001
$ws.Range("L122").Value = 6513
$ws.Range("M122").Value = -1793.4001
$ws.Range("N122").Value = -11413
$ws.Range("H132").Value = 2553.8147
$ws.Range("I132").Value = 2137.15
$ws.Range("J132").Value = 3744.2856
$ws.Range("K132").Value = 6411.450000000001
$ws.Range("L132").Value = 11232.8568
$ws.Range("M132").Value = -3881.450000000001
$ws.Range("N132").Value = -16292.8568

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7598.7334
$ws.Range("I3").Value = 912.9167
$ws.Range("J3").Value = 34342
$ws.Range("K3").Value = 912.9167
$ws.Range("L3").Value = 34342
$ws.Range("M3").Value = -798.9167
$ws.Range("N3").Value = -34570
$ws.Range("H100").Value = 10000
$ws.Range("J100").Value = 10000
$ws.Range("L100").Value = 10000
$ws.Range("N100").Value = -12164
$ws.Range("H106").Value = 25780.334
$ws.Range("J106").Value = 25780.334
$ws.Range("L106").Value = 25780.334
$ws.Range("N106").Value = -28304.334
$ws.Range("H107").Value = 1887.35
$ws.Range("I107").Value = 1357.1
$ws.Range("K107").Value = 1357.1
$ws.Range("M107").Value = 562.9000000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 13980.5
$ws.Range("I93").Value = 7143.857
$ws.Range("J93").Value = 29932.666
$ws.Range("K93").Value = 7143.857
$ws.Range("L93").Value = 29932.666
$ws.Range("M93").Value = -5271.857
$ws.Range("N93").Value = -33676.666
$ws.Range("H99").Value = 1484.25
$ws.Range("I99").Value = 1554.8889
$ws.Range("J99").Value = 1393.4286
$ws.Range("K99").Value = 1554.8889
$ws.Range("L99").Value = 1393.4286
$ws.Range("M99").Value = -56.88889999999992
$ws.Range("N99").Value = -4389.4286
$ws.Range("H122").Value = 857.1429000000001
$ws.Range("I122").Value = 700
$ws.Range("K122").Value = 2100
$ws.Range("M122").Value = 350
$ws.Range("H126").Value = 1484.25
$ws.Range("I126").Value = 1554.8889
$ws.Range("J126").Value = 1393.4286
$ws.Range("K126").Value = 4664.6667
$ws.Range("L126").Value = 4180.2858
$ws.Range("M126").Value = -2194.6667
$ws.Range("N126").Value = -9120.2858
$ws.Range("H132").Value = 7033.85
$ws.Range("I132").Value = 10225.818
$ws.Range("J132").Value = 3132.5557
$ws.Range("K132").Value = 30677.454
$ws.Range("L132").Value = 9397.667099999999
$ws.Range("M132").Value = -28147.454
$ws.Range("N132").Value = -14457.6671

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 669
$ws.Range("I14").Value = 669
$ws.Range("K14").Value = 2007
$ws.Range("M14").Value = -1834

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 165.2
$ws.Range("I2").Value = 91.666664
$ws.Range("J2").Value = 275.5
$ws.Range("K2").Value = 91.666664
$ws.Range("L2").Value = 275.5
$ws.Range("M2").Value = 21.333336
$ws.Range("N2").Value = -501.5
$ws.Range("H97").Value = 500
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 500
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1492
$ws.Range("H113").Value = 2161.2083
$ws.Range("I113").Value = 1319.2222
$ws.Range("J113").Value = 2666.4
$ws.Range("K113").Value = 1319.2222
$ws.Range("L113").Value = 2666.4
$ws.Range("M113").Value = 850.7778000000001
$ws.Range("N113").Value = -7006.4

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1719.875
$ws.Range("J22").Value = 1939.75
$ws.Range("L22").Value = 1939.75
$ws.Range("N22").Value = -2529.75
$ws.Range("H27").Value = 1719.875
$ws.Range("J27").Value = 1939.75
$ws.Range("L27").Value = 1939.75
$ws.Range("N27").Value = -2153.75
$ws.Range("H40").Value = 3227.1667
$ws.Range("I40").Value = 2940.75
$ws.Range("K40").Value = 2940.75
$ws.Range("M40").Value = -2804.75
$ws.Range("H105").Value = 19750
$ws.Range("J105").Value = 19750
$ws.Range("L105").Value = 19750
$ws.Range("N105").Value = -26738
$ws.Range("H136").Value = 7133.9443
$ws.Range("I136").Value = 8931.615
$ws.Range("K136").Value = 26794.845
$ws.Range("M136").Value = -24244.845

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 7124.75
$ws.Range("J101").Value = 7124.75
$ws.Range("L101").Value = 7124.75
$ws.Range("N101").Value = -13614.75
$ws.Range("H133").Value = 29000
$ws.Range("J133").Value = 29000
$ws.Range("L133").Value = 29000
$ws.Range("N133").Value = -39120
